$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2098.125
$ws.Range("I40").Value = 2043.7778
$ws.Range("J40").Value = 2168
$ws.Range("K40").Value = 2043.7778
$ws.Range("L40").Value = 2168
$ws.Range("M40").Value = -1868.7778
$ws.Range("N40").Value = -2518
$ws.Range("H53").Value = 256.0909
$ws.Range("I53").Value = 226
$ws.Range("K53").Value = 226
$ws.Range("M53").Value = 411
$ws.Range("H111").Value = 2353
$ws.Range("I111").Value = 442.66666
$ws.Range("K111").Value = 1327.99998
$ws.Range("M111").Value = 1739.00002
$ws.Range("H112").Value = 3265.6667
$ws.Range("J112").Value = 3999
$ws.Range("L112").Value = 11997
$ws.Range("N112").Value = -14213
$ws.Range("H113").Value = 2837.5
$ws.Range("I113").Value = 2487.875
$ws.Range("K113").Value = 2487.875
$ws.Range("M113").Value = 766.125
$ws.Range("H132").Value = 1207.6923
$ws.Range("I132").Value = 1155.091
$ws.Range("K132").Value = 3465.273
$ws.Range("M132").Value = -935.2729999999997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 580
$ws.Range("I2").Value = 388
$ws.Range("K2").Value = 388
$ws.Range("M2").Value = -275
$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 500
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 500
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -384
$ws.Range("N4").Value = ""
$ws.Range("H9").Value = 2000000
$ws.Range("I9").Value = 2000000
$ws.Range("K9").Value = 2000000
$ws.Range("M9").Value = -1999830
$ws.Range("H20").Value = 2000000
$ws.Range("I20").Value = 2000000
$ws.Range("K20").Value = 2000000
$ws.Range("M20").Value = -1999730
$ws.Range("H63").Value = 6668.3335
$ws.Range("I63").Value = 6668.3335
$ws.Range("K63").Value = 6668.3335
$ws.Range("M63").Value = -5982.3335
$ws.Range("H66").Value = 6668.3335
$ws.Range("I66").Value = 6668.3335
$ws.Range("K66").Value = 33341.6675
$ws.Range("M66").Value = -29909.6675
$ws.Range("H74").Value = 2224.9546
$ws.Range("I74").Value = 1804.75
$ws.Range("K74").Value = 1804.75
$ws.Range("M74").Value = -930.75
$ws.Range("H77").Value = 2224.9546
$ws.Range("I77").Value = 1804.75
$ws.Range("K77").Value = 9023.75
$ws.Range("M77").Value = -4655.75
$ws.Range("H116").Value = 580
$ws.Range("I116").Value = 388
$ws.Range("K116").Value = 388
$ws.Range("M116").Value = 1906
$ws.Range("H132").Value = 1473.9
$ws.Range("I132").Value = 1493.7241
$ws.Range("K132").Value = 4481.1723
$ws.Range("M132").Value = -1951.1723

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 580
$ws.Range("I3").Value = 388
$ws.Range("K3").Value = 388
$ws.Range("M3").Value = -274
$ws.Range("H20").Value = 696.6667
$ws.Range("I20").Value = 697.5
$ws.Range("J20").Value = 695
$ws.Range("K20").Value = 697.5
$ws.Range("L20").Value = 695
$ws.Range("M20").Value = -450.5
$ws.Range("N20").Value = -1189
$ws.Range("H22").Value = 1025.6
$ws.Range("I22").Value = 796
$ws.Range("J22").Value = 1370
$ws.Range("K22").Value = 796
$ws.Range("L22").Value = 1370
$ws.Range("M22").Value = -623
$ws.Range("N22").Value = -1716
$ws.Range("H99").Value = 1574.25
$ws.Range("I99").Value = 1098.375
$ws.Range("J99").Value = 3477.75
$ws.Range("K99").Value = 1098.375
$ws.Range("L99").Value = 3477.75
$ws.Range("M99").Value = 399.625
$ws.Range("N99").Value = -6473.75
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = ""
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 134.03847
$ws.Range("I7").Value = 124.22222
$ws.Range("J7").Value = 156.125
$ws.Range("K7").Value = 124.22222
$ws.Range("L7").Value = 156.125
$ws.Range("M7").Value = -11.22221999999999
$ws.Range("N7").Value = -382.125
$ws.Range("H22").Value = 31891.375
$ws.Range("J22").Value = 50508.6
$ws.Range("L22").Value = 50508.6
$ws.Range("N22").Value = -51208.6
$ws.Range("H31").Value = 1752.5883
$ws.Range("J31").Value = 4220.5
$ws.Range("L31").Value = 4220.5
$ws.Range("N31").Value = -4810.5
$ws.Range("H34").Value = 1752.5883
$ws.Range("J34").Value = 4220.5
$ws.Range("L34").Value = 4220.5
$ws.Range("N34").Value = -4624.5
$ws.Range("H55").Value = 13333.333
$ws.Range("I55").Value = 10000
$ws.Range("J55").Value = 15000
$ws.Range("K55").Value = 10000
$ws.Range("L55").Value = 15000
$ws.Range("M55").Value = -9685
$ws.Range("N55").Value = -15630
$ws.Range("H94").Value = 126915.89
$ws.Range("I94").Value = 224006.4
$ws.Range("K94").Value = 224006.4
$ws.Range("M94").Value = -223555.4
$ws.Range("H99").Value = 2733.8
$ws.Range("I99").Value = 1184.3334
$ws.Range("J99").Value = 3397.8572
$ws.Range("K99").Value = 1184.3334
$ws.Range("L99").Value = 3397.8572
$ws.Range("M99").Value = 313.6666
$ws.Range("N99").Value = -6393.8572
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = ""
$ws.Range("H126").Value = 2733.8
$ws.Range("I126").Value = 1184.3334
$ws.Range("J126").Value = 3397.8572
$ws.Range("K126").Value = 3553.0002
$ws.Range("L126").Value = 10193.5716
$ws.Range("M126").Value = -1083.0002
$ws.Range("N126").Value = -15133.5716
$ws.Range("H132").Value = 1386.579
$ws.Range("I132").Value = 1262.4062
$ws.Range("J132").Value = 2048.8333
$ws.Range("K132").Value = 3787.2186
$ws.Range("L132").Value = 6146.499899999999
$ws.Range("M132").Value = -1257.2186
$ws.Range("N132").Value = -11206.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 2333.6667
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").Value = ""
$ws.Range("H86").Value = 4798.1665
$ws.Range("I86").Value = 929.6667
$ws.Range("J86").Value = 8666.666999999999
$ws.Range("K86").Value = 2789.0001
$ws.Range("L86").Value = 26000.001
$ws.Range("M86").Value = -1603.0001
$ws.Range("N86").Value = -28372.001
$ws.Range("H89").Value = 4798.1665
$ws.Range("I89").Value = 929.6667
$ws.Range("J89").Value = 8666.666999999999
$ws.Range("K89").Value = 8367.0003
$ws.Range("L89").Value = 78000.003
$ws.Range("M89").Value = -2439.0003
$ws.Range("N89").Value = -89856.003
$ws.Range("H139").Value = 2271.4
$ws.Range("I139").Value = 2271.4
$ws.Range("K139").Value = 6814.200000000001
$ws.Range("M139").Value = -1674.200000000001
$ws.Range("H140").Value = 6249.4116
$ws.Range("I140").Value = 1476.6364
$ws.Range("K140").Value = 4429.9092
$ws.Range("M140").Value = 750.0907999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 26000
$ws.Range("J105").Value = 26000
$ws.Range("L105").Value = 26000
$ws.Range("N105").Value = -32988
$ws.Range("H113").Value = 1276.6666
$ws.Range("I113").Value = 1276.6666
$ws.Range("K113").Value = 1276.6666
$ws.Range("M113").Value = 893.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1516.6666
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = ""
$ws.Range("H27").Value = 1516.6666
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").Value = ""
$ws.Range("H40").Value = 4145.8335
$ws.Range("H61").Value = 2665.9167
$ws.Range("I61").Value = 2322
$ws.Range("K61").Value = 2322
$ws.Range("M61").Value = -2120
$ws.Range("H96").Value = 50449.5
$ws.Range("J96").Value = 50449.5
$ws.Range("L96").Value = 50449.5
$ws.Range("N96").Value = -55941.5
$ws.Range("H113").Value = 2665.9167
$ws.Range("I113").Value = 2322
$ws.Range("K113").Value = 2322
$ws.Range("M113").Value = -152
$ws.Range("H122").Value = 7842.7617
$ws.Range("I122").Value = 7907.2856
$ws.Range("K122").Value = 23721.8568
$ws.Range("M122").Value = -21271.8568
$ws.Range("H132").Value = 3044.6086
$ws.Range("I132").Value = 2592.4119
$ws.Range("J132").Value = 4325.8335
$ws.Range("K132").Value = 7777.2357
$ws.Range("L132").Value = 12977.5005
$ws.Range("M132").Value = -5247.2357
$ws.Range("N132").Value = -18037.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 11114431
$ws.Range("I100").Value = 12502484
$ws.Range("K100").Value = 25004968
$ws.Range("M100").Value = -25004427
$ws.Range("H122").Value = 345.2
$ws.Range("I122").Value = 345.2
$ws.Range("K122").Value = 1035.6
$ws.Range("M122").Value = 1414.4
